$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70: Consecrating Congregation | Holy Water
$ws.Range("H70").Value = 2751.125
$ws.Range("I70").Value = 2900.3333
$ws.Range("J70").Value = 2661.6
$ws.Range("K70").Value = 8700.999899999999
$ws.Range("L70").Value = 7984.799999999999
$ws.Range("M70").Value = -8430.999899999999
$ws.Range("N70").Value = -8524.799999999999

# Row 73: Curbing the Contagion (L) | Holy Water
$ws.Range("H73").Value = 2751.125
$ws.Range("I73").Value = 2900.3333
$ws.Range("J73").Value = 2661.6
$ws.Range("K73").Value = 8700.999899999999
$ws.Range("L73").Value = 7984.799999999999
$ws.Range("M73").Value = -7764.999899999999
$ws.Range("N73").Value = -9856.799999999999

# Row 86: Filling in the Blanks | Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 3105.6667
$ws.Range("I86").Value = 3176.8
$ws.Range("J86").Value = 2750
$ws.Range("K86").Value = 3176.8
$ws.Range("L86").Value = 2750
$ws.Range("M86").Value = -2053.8
$ws.Range("N86").Value = -4996

# Row 88: The Grave of Hemlock Groves | Growth Formula Zeta
$ws.Range("H88").Value = 11509.818
$ws.Range("I88").Value = 3000
$ws.Range("J88").Value = 12360.8
$ws.Range("K88").Value = 3000
$ws.Range("L88").Value = 12360.8
$ws.Range("M88").Value = -2594
$ws.Range("N88").Value = -13172.8

# Row 89: Ink into Antiquity (L) | Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 3105.6667
$ws.Range("I89").Value = 3176.8
$ws.Range("J89").Value = 2750
$ws.Range("K89").Value = 15884
$ws.Range("L89").Value = 13750
$ws.Range("M89").Value = -10268
$ws.Range("N89").Value = -24982

# Row 91: Dappling the Highlands (L) | Growth Formula Zeta
$ws.Range("H91").Value = 11509.818
$ws.Range("I91").Value = 3000
$ws.Range("J91").Value = 12360.8
$ws.Range("K91").Value = 3000
$ws.Range("L91").Value = 12360.8
$ws.Range("M91").Value = -1596
$ws.Range("N91").Value = -15168.8

# Row 129: Practical Command | Commanding Craftsman's Draught
$ws.Range("H129").Value = 1146.6046
$ws.Range("I129").Value = 816.25
$ws.Range("J129").Value = 1222.1143
$ws.Range("K129").Value = 2448.75
$ws.Range("L129").Value = 3666.3429
$ws.Range("M129").Value = 2551.25
$ws.Range("N129").Value = -13666.3429

# Row 131: Mindful Study | Grade 5 Tincture of Mind
$ws.Range("H131").Value = 1770.5
$ws.Range("I131").Value = 491.25
$ws.Range("J131").Value = 3049.75
$ws.Range("K131").Value = 1473.75
$ws.Range("L131").Value = 9149.25
$ws.Range("M131").Value = 3566.25
$ws.Range("N131").Value = -19229.25

$ws = $wb.Worksheets.Item("ARM")
# Row 29: No Hand-me-downs | Iron Vambraces
$ws.Range("H29").Value = 8140
$ws.Range("J29").Value = 8140
$ws.Range("L29").Value = 8140
$ws.Range("N29").Value = -8756

# Row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 1071527.4
$ws.Range("I122").Value = 1427953.1
$ws.Range("J122").Value = 2250
$ws.Range("K122").Value = 4283859.300000001
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -4281409.300000001
$ws.Range("N122").Value = -11650

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 3108.8276
$ws.Range("I132").Value = 2326.4243
$ws.Range("J132").Value = 4141.6
$ws.Range("K132").Value = 6979.2729
$ws.Range("L132").Value = 12424.8
$ws.Range("M132").Value = -4449.2729
$ws.Range("N132").Value = -17484.8

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 2016.4736
$ws.Range("I86").Value = 1892.1666
$ws.Range("J86").Value = 2229.5715
$ws.Range("K86").Value = 1892.1666
$ws.Range("L86").Value = 2229.5715
$ws.Range("M86").Value = -769.1666
$ws.Range("N86").Value = -4475.5715

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 2016.4736
$ws.Range("I89").Value = 1892.1666
$ws.Range("J89").Value = 2229.5715
$ws.Range("K89").Value = 9460.833000000001
$ws.Range("L89").Value = 11147.8575
$ws.Range("M89").Value = -3844.833000000001
$ws.Range("N89").Value = -22379.8575

$ws = $wb.Worksheets.Item("CUL")
# Row 21: Shy Is the Oyster | Raw Oyster
$ws.Range("H21").Value = 969.9
$ws.Range("I21").Value = 300
$ws.Range("J21").Value = 1416.5
$ws.Range("K21").Value = 900
$ws.Range("L21").Value = 4249.5
$ws.Range("M21").Value = -727
$ws.Range("N21").Value = -4595.5

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 1563478.2
$ws.Range("J131").Value = 1064.7858
$ws.Range("L131").Value = 3194.3574
$ws.Range("N131").Value = -13274.3574

# Row 140: Sweet, Sweet Bean Juice | Mesquite Juice
$ws.Range("H140").Value = 2375.8333
$ws.Range("I140").Value = 2474.7368
$ws.Range("J140").Value = 2000
$ws.Range("K140").Value = 7424.2104
$ws.Range("L140").Value = 6000
$ws.Range("M140").Value = -2244.2104
$ws.Range("N140").Value = -16360

$ws = $wb.Worksheets.Item("GSM")
# Row 63: Not on My Table | Mythrite Earrings of Healing
$ws.Range("H63").Value = 40000
$ws.Range("J63").Value = 40000
$ws.Range("L63").Value = 40000
$ws.Range("N63").Value = -41372

# Row 66: Heinz's Dilemma (L) | Mythrite Earrings of Healing
$ws.Range("H66").Value = 40000
$ws.Range("J66").Value = 40000
$ws.Range("L66").Value = 120000
$ws.Range("N66").Value = -126864

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 3203.7805
$ws.Range("I132").Value = 3748.3333
$ws.Range("J132").Value = 2632
$ws.Range("K132").Value = 11244.9999
$ws.Range("L132").Value = 7896
$ws.Range("M132").Value = -8714.999899999999
$ws.Range("N132").Value = -12956

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs | Aldgoat Leather
$ws.Range("H22").Value = 1994.5186
$ws.Range("I22").Value = 558.6667
$ws.Range("J22").Value = 2404.762
$ws.Range("K22").Value = 558.6667
$ws.Range("L22").Value = 2404.762
$ws.Range("M22").Value = -263.6667
$ws.Range("N22").Value = -2994.762

# Row 27: Fire and Hide | Aldgoat Leather
$ws.Range("H27").Value = 1994.5186
$ws.Range("I27").Value = 558.6667
$ws.Range("J27").Value = 2404.762
$ws.Range("K27").Value = 558.6667
$ws.Range("L27").Value = 2404.762
$ws.Range("M27").Value = -451.6667
$ws.Range("N27").Value = -2618.762

# Row 40: Best Served Toad | Toad Leather
$ws.Range("H40").Value = 25001854
$ws.Range("I40").Value = 31251870
$ws.Range("J40").Value = 1788.125
$ws.Range("K40").Value = 31251870
$ws.Range("L40").Value = 1788.125
$ws.Range("M40").Value = -31251734
$ws.Range("N40").Value = -2060.125

# Row 46: Supply Side Logic | Boar Leather
$ws.Range("H46").Value = 1294.9
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1294.9
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1294.9
$ws.Range("M46").Value = $null
$ws.Range("N46").Value = -1670.9

# Row 61: Spelling Me Softly | Raptor Leather
$ws.Range("H61").Value = 1513.9584
$ws.Range("J61").Value = 1822
$ws.Range("L61").Value = 1822
$ws.Range("N61").Value = -2226

# Row 113: Peace in Rest | Atrociraptor Leather
$ws.Range("H113").Value = 1513.9584
$ws.Range("J113").Value = 1822
$ws.Range("L113").Value = 1822
$ws.Range("N113").Value = -6162

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 13339760
$ws.Range("J132").Value = 2357.4
$ws.Range("L132").Value = 7072.200000000001
$ws.Range("N132").Value = -12132.2

$ws = $wb.Worksheets.Item("WVR")
# Row 64: Ribbon of Remembrance | Rainbow Ribbon of Healing
$ws.Range("H64").Value = 31055
$ws.Range("J64").Value = 31055
$ws.Range("L64").Value = 31055
$ws.Range("N64").Value = -31551

# Row 67: The Road Was a Ribbon of Moonlight (L) | Rainbow Ribbon of Healing
$ws.Range("H67").Value = 31055
$ws.Range("J67").Value = 31055
$ws.Range("L67").Value = 31055
$ws.Range("N67").Value = -32771

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 2637.125
$ws.Range("I132").Value = 1927.3636
$ws.Range("J132").Value = 4198.6
$ws.Range("K132").Value = 5782.0908
$ws.Range("L132").Value = 12595.8
$ws.Range("M132").Value = -3252.0908
$ws.Range("N132").Value = -17655.8
